$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.267898082733154
$ws.Range("B1").Value = 2.809709072113037
$ws.Range("C1").Value = 8.678452491760254
$ws.Range("D1").Value = 2.01676082611084
$ws.Range("E1").Value = 1.128594517707825
